$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.309.48"
$ws.Range("E2").Value = "  +2.41%  "
$ws.Range("D3").Value = "3.389.96"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'586.20"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").Value = "'180.69"
$ws.Range("E6").Value = "  +2.86%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("D9").Value = "'0.198"
$ws.Range("E9").Value = "  +8.16%  "
$ws.Range("E10").Value = "  +2.36%  "
$ws.Range("D11").Value = "'48.65"
$ws.Range("E11").Value = "  +3.80%  "
$ws.Range("D12").Value = "'0.0000284"
$ws.Range("E12").Value = "  +4.32%  "
$ws.Range("D13").Value = "'679.69"
$ws.Range("E13").Value = "  -1.93%  "
$ws.Range("D14").Value = "'8.66"
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("D15").Value = "3.933.03"
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("D16").Value = "69.353.50"
$ws.Range("D17").Value = "3.404.75"
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("E18").Value = "  +1.56%  "
$ws.Range("D19").Value = "'17.74"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("E20").Value = "  +2.33%  "
$ws.Range("D21").Value = "'0.907"
$ws.Range("E21").Value = "  +1.50%  "
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").Value = "'17.18"
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("D24").Value = "'103.06"
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("D25").Value = "'3.93"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").Value = "'9.64"
$ws.Range("E27").Value = "  +2.12%  "
$ws.Range("D28").Value = "'33.91"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("D29").Value = "'8.78"
$ws.Range("E29").Value = "  +2.74%  "
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").Value = "'11.16"
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("D32").Value = "'556.36"
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("E33").Value = "  +10.35%  "
$ws.Range("E34").Value = "  +1.32%  "
$ws.Range("D35").Value = "'58.55"
$ws.Range("E35").Value = "  +2.43%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").Value = "3.670.60"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("D38").Value = "'0.139"
$ws.Range("E38").Value = "  +4.55%  "
$ws.Range("D39").Value = "'35.59"
$ws.Range("E39").Value = "  +1.19%  "
$ws.Range("D40").Value = "0.0₃0721"
$ws.Range("E40").Value = "  +7.41%  "
$ws.Range("D41").Value = "'3.28"
$ws.Range("E41").Value = "  +3.62%  "
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("D43").Value = "'0.340"
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("E44").Value = "  +3.88%  "
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("E48").Value = "  +5.50%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").Value = "'133.64"
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("E51").Value = "  +4.83%  "
